$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Updated fitting parameters
$ws.Range("J2").Value = 0.034470000000000001
$ws.Range("K2").Value = 0.017500000000000002

# Mark the "-" unit cell for K (o_w) as ready / underline it, like the row above
$ws.Range("K3").Font.Underline = $true
